# Update TPM-derived NATMI metrics for Gnai2-Cxcr2 (rows 2-13)
# Values below are taken verbatim from the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 0.01989833333333333
$ws.Range("N2").Value = 0.059695
$ws.Range("O2").Value = 0.5455534129646046
$ws.Range("P2").Value = 0.5455534129646046
$ws.Range("Q2").Value = 3.381038381538334
$ws.Range("R2").Value = 30.429345433845
$ws.Range("S2").Value = 0.2422903597114736
$ws.Range("T2").Value = 0.2422903597114736

# Row 3
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.1997148627777118
$ws.Range("P3").Value = 0.1997148627777118
$ws.Range("Q3").Value = 1.237722284140333
$ws.Range("R3").Value = 11.139500557263
$ws.Range("S3").Value = 0.08869706392117986
$ws.Range("T3").Value = 0.08869706392117989

# Row 4
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 0.009290999999999999
$ws.Range("N4").Value = 0.027873
$ws.Range("O4").Value = 0.2547317242576836
$ws.Range("P4").Value = 0.2547317242576836
$ws.Range("Q4").Value = 1.578686369187
$ws.Range("R4").Value = 14.208177322683
$ws.Range("S4").Value = 0.1131310695407974
$ws.Range("T4").Value = 0.1131310695407975

# Row 5
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("M5").Value = 0.01989833333333333
$ws.Range("N5").Value = 0.059695
$ws.Range("O5").Value = 0.5455534129646046
$ws.Range("P5").Value = 0.5455534129646046
$ws.Range("Q5").Value = 1.360692665295
$ws.Range("R5").Value = 12.246233987655
$ws.Range("S5").Value = 0.09750930871748562
$ws.Range("T5").Value = 0.09750930871748562

# Row 6
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("O6").Value = 0.1997148627777118
$ws.Range("P6").Value = 0.1997148627777118
$ws.Range("S6").Value = 0.03569596990373086
$ws.Range("T6").Value = 0.03569596990373086

# Row 7
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("M7").Value = 0.009290999999999999
$ws.Range("N7").Value = 0.027873
$ws.Range("O7").Value = 0.2547317242576836
$ws.Range("P7").Value = 0.2547317242576836
$ws.Range("Q7").Value = 0.6353394197129999
$ws.Range("R7").Value = 5.718054777417
$ws.Range("S7").Value = 0.04552939043274104
$ws.Range("T7").Value = 0.04552939043274105

# Row 8
$ws.Range("G8").Value = 53.27463399999999
$ws.Range("H8").Value = 159.823902
$ws.Range("I8").Value = 0.1392470275793777
$ws.Range("J8").Value = 0.1392470275793778
$ws.Range("M8").Value = 0.01989833333333333
$ws.Range("N8").Value = 0.059695
$ws.Range("O8").Value = 0.5455534129646046
$ws.Range("P8").Value = 0.5455534129646046
$ws.Range("Q8").Value = 1.060076425543333
$ws.Range("R8").Value = 9.540687829889999
$ws.Range("S8").Value = 0.07596669114110595
$ws.Range("T8").Value = 0.07596669114110596

# Row 9
$ws.Range("G9").Value = 53.27463399999999
$ws.Range("H9").Value = 159.823902
$ws.Range("I9").Value = 0.1392470275793777
$ws.Range("J9").Value = 0.1392470275793778
$ws.Range("O9").Value = 0.1997148627777118
$ws.Range("P9").Value = 0.1997148627777118
$ws.Range("Q9").Value = 0.3880701922673333
$ws.Range("R9").Value = 3.492631730406
$ws.Range("S9").Value = 0.02780970100521967
$ws.Range("T9").Value = 0.02780970100521968

# Row 10
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 0.009290999999999999
$ws.Range("N10").Value = 0.027873
$ws.Range("O10").Value = 0.2547317242576836
$ws.Range("P10").Value = 0.2547317242576836
$ws.Range("Q10").Value = 0.4949746244939999
$ws.Range("R10").Value = 4.454771620446
$ws.Range("S10").Value = 0.03547063543305211
$ws.Range("T10").Value = 0.03547063543305213

# Row 11
$ws.Range("G11").Value = 91.01828266666666
$ws.Range("H11").Value = 273.054848
$ws.Range("I11").Value = 0.2378998101932138
$ws.Range("J11").Value = 0.2378998101932138
$ws.Range("M11").Value = 0.01989833333333333
$ws.Range("N11").Value = 0.059695
$ws.Range("O11").Value = 0.5455534129646046
$ws.Range("P11").Value = 0.5455534129646046
$ws.Range("Q11").Value = 1.811112127928889
$ws.Range("R11").Value = 16.30000915136
$ws.Range("S11").Value = 0.1297870533945394
$ws.Range("T11").Value = 0.1297870533945394

# Row 12
$ws.Range("G12").Value = 91.01828266666666
$ws.Range("H12").Value = 273.054848
$ws.Range("I12").Value = 0.2378998101932138
$ws.Range("J12").Value = 0.2378998101932138
$ws.Range("O12").Value = 0.1997148627777118
$ws.Range("P12").Value = 0.1997148627777118
$ws.Range("Q12").Value = 0.6630075103715556
$ws.Range("R12").Value = 5.967067593344
$ws.Range("S12").Value = 0.04751212794758137
$ws.Range("T12").Value = 0.04751212794758138

# Row 13
$ws.Range("G13").Value = 91.01828266666666
$ws.Range("H13").Value = 273.054848
$ws.Range("I13").Value = 0.2378998101932138
$ws.Range("J13").Value = 0.2378998101932138
$ws.Range("M13").Value = 0.009290999999999999
$ws.Range("N13").Value = 0.027873
$ws.Range("O13").Value = 0.2547317242576836
$ws.Range("P13").Value = 0.2547317242576836
$ws.Range("Q13").Value = 0.8456508642559999
$ws.Range("R13").Value = 7.610857778303999
$ws.Range("S13").Value = 0.060600628851093
$ws.Range("T13").Value = 0.06060062885109302
